$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldName = "Administrator, Miss Dina Nasr"
$newName = "Miss Dina Nasr, Administrator"

$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq $oldName) {
        $cell.Value = $newName
    }
}
